# Apply updated cryptocurrency price/volume data to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '90.775.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.116.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.14%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.15%  '

$ws.Range('E7').Value = '  +11.43%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.372'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.11%  '

$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('B10').Value = 'LidoStakedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.113.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.95%  '

$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.752'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.27%  '

$ws.Range('E12').Value = '  +3.79%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.54%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.61%  '

$ws.Range('E15').Value = '  -1.90%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.542.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.684.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.117.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.82%  '

$ws.Range('E19').Value = '  +4.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.61%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000211'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.95%  '

$ws.Range('E22').Value = '  +7.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '448.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.25%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.00%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.68%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '93.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.262.29'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('E30').Value = '  +11.87%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.220'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +12.71%  '

$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +35.36%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +27.11%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.159'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.65%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.00%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.66'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.46%  '

$ws.Range('B38').Value = 'MantraDAO'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.22'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +29.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '494.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.61%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.17%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.418'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.26%  '

$ws.Range('E45').Value = '  -0.05%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.687'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.38%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.97%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.09%  '

$ws.Range('E51').Value = '  -0.26%  '
